# Word COM-interop script implementing the frontmatter.docx "List of Figures"
# renumbering described in the commit message: "reworking 7.1 and all
# references to it. moving peripheral settings details to appendices.
# updating website"
#
# Net effect on the "Figures in Chapter 3" list: every caption from
# Figure 3.3 onward shifts up into the previous figure's slot, Figure 3.2's
# old caption is discarded, Figure 3.9 is reworded along the way ("Ideation
# Grids" -> "Ideation Decks" when it moves to slot 3.8), and Figure 3.14
# receives a brand-new caption.
#
# We process the paragraphs from the LAST figure to the FIRST so that the
# text we search for with Find is always still the original, untouched
# text (avoiding collisions with text we've already written earlier in the
# list).

$d = $word.ActiveDocument

$dash = [char]0x2013      # en dash
$lq   = [char]0x201C      # left curly quote
$rq   = [char]0x201D      # right curly quote

function Set-CaptionAfterHyperlink {
    param($OldText, $NewPieces)

    $r = $d.Content
    $found = $r.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find caption text: $OldText"
    }

    $pos = $r.Start
    # Clear the matched text completely first so the following inserts
    # start from a clean slate (this avoids leftover/merged runs).
    $r.Text = ""

    foreach ($piece in $NewPieces) {
        $ins = $d.Range($pos, $pos)
        $ins.InsertAfter($piece)
        $pos = $pos + $piece.Length
    }
}

# Figure 3.14 (was "Pilot Study Recruitment Poster")
Set-CaptionAfterHyperlink "- Pilot Study Recruitment Poster" `
    @("- How the Case Studies and Peripheral Activities Contribute to This Thesis")

# Figure 3.13 (was "Spreadsheet-based Quantitative Analysis of Interview Data for Case Study Two")
Set-CaptionAfterHyperlink "- Spreadsheet-based Quantitative Analysis of Interview Data for Case Study Two" `
    @("- Pilot Study Recruitment Poster")

# Figure 3.12 (was "Thematic Analysis of Qualitative Data using Quirkos for Case Study One")
Set-CaptionAfterHyperlink "- Thematic Analysis of Qualitative Data using Quirkos for Case Study One" `
    @("- Spreadsheet-based Quantitative Analysis of Interview Data for Case Study Two")

# Figure 3.11 (was the Storyboarding Cards caption)
Set-CaptionAfterHyperlink ("- Storyboarding Cards " + $dash + " A Collaboratively-constructed Narrative Created through Discussion From a Palette of Possible Parent and Staff Actions") `
    @("- Thematic Analysis of Qualitative Data using Quirkos for Case Study One")

# Figure 3.10 (was the Group Poster Design caption)
Set-CaptionAfterHyperlink ("- Group Poster Design " + $dash + " A Participant-designed Poster to Advertise Features of Imagined Data Interface Products") `
    @(("- Storyboarding Cards " + $dash + " A Collaboratively-constructed Narrative Created through Discussion From a Palette of Possible Parent and Staff Actions"))

# Figure 3.9 (was "Ideation Grids ...")
Set-CaptionAfterHyperlink ("- Ideation Grids " + $dash + " Combining Random Design Ingredients to Generate New Ideas") `
    @(("- Group Poster Design " + $dash + " A Participant-designed Poster to Advertise Features of Imagined Data Interface Products"))

# Figure 3.8 (was the Home Interviewing caption) -> note wording change Grids -> Decks
Set-CaptionAfterHyperlink "- Home Interviewing: Card Sorting With a Family in Their Living Room" `
    @(("- Ideation Decks " + $dash + " Combining Random Design Ingredients to Generate New Ideas"))

# Figure 3.7 (was the Personal Data Examples caption)
Set-CaptionAfterHyperlink ("- Personal Data Examples " + $dash + " Making Data Relatable") `
    @("- Home Interviewing: Card Sorting With a Family in Their Living Room")

# Figure 3.6 (was the Family Civic Data Cards caption)
Set-CaptionAfterHyperlink ("- Family Civic Data Cards " + $dash + " Things to Think With") `
    @(("- Personal Data Examples " + $dash + " Making Data Relatable"))

# Figure 3.5 (was the Sentence Ranking caption)
Set-CaptionAfterHyperlink ("- Sentence Ranking " + $dash + " Bringing Support Workers and Families to a Shared Problem Space") `
    @(("- Family Civic Data Cards " + $dash + " Things to Think With"))

# Figure 3.4 (was the Walls of Data caption)
Set-CaptionAfterHyperlink ("- Walls of Data " + $dash + " Sensitising Participants to the World of Commercially-held Data and GDPR") `
    @(("- Sentence Ranking " + $dash + " Bringing Support Workers and Families to a Shared Problem Space"))

# Figure 3.3 (was the multi-run "- <U+201C>Family Facts<U+201D> - What is Data?" caption)
Set-CaptionAfterHyperlink ("- " + $lq + "Family Facts" + $rq + " " + $dash + " What is Data?") `
    @(("- Walls of Data " + $dash + " Sensitising Participants to the World of Commercially-held Data and GDPR"))

# Figure 3.2 (was "Research Activities and Contexts") -> becomes the
# multi-run "- <U+201C>Family Facts<U+201D> - What is Data?" caption (split into the
# same seven runs the original Figure 3.3 paragraph used).
Set-CaptionAfterHyperlink "- Research Activities and Contexts" `
    @("-", " ", $lq, "Family Facts", $rq, " ", ($dash + " What is Data?"))
